# Apply the edit described by the diff:
# - Header A1: "File Name" -> "Loc"
# - Header B1: "Unnormalized P_max" -> "P_max"
# - Header C1: "Electrode Locations" removed (column C deleted)
# - For each data row, column A (the long filename) is replaced by the
#   electrode location that used to live in column C, and column C is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row from the sheet's UsedRange.
$lastRow = $ws.UsedRange.Rows.Count

# Replace column A values (rows 2..lastRow) with the electrode-location values
# that are currently stored in column C, before column C gets removed.
for ($r = 2; $r -le $lastRow; $r++) {
    $locCell = $ws.Cells.Item($r, 3)
    $loc = $locCell.Value()
    $ws.Cells.Item($r, 1).Value = $loc
}

# Update header text.
$ws.Cells.Item(1, 1).Value = "Loc"
$ws.Cells.Item(1, 2).Value = "P_max"

# Remove column C entirely (also removes the "Electrode Locations" header).
$ws.Columns.Item(3).Delete()
